$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44873
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 22000
$ws.Range("O2").Value = 22500
$ws.Range("P2").Value = 22250
$ws.Range("S2").Value = 2781
$ws.Range("D3").Value = 44495
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 270
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19556
$ws.Range("S3").Value = 2444
$ws.Range("D4").Value = 44523
$ws.Range("M4").Value = 400
$ws.Range("P4").Value = 21500
$ws.Range("S4").Value = 2688
$ws.Range("D5").Value = 44523
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("Q5").Value = "$/bandeja 8 kilos"
$ws.Range("S5").Value = 2250
$ws.Range("T5").Value = 8
$ws.Range("D6").Value = 44519
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 400
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 2688
$ws.Range("D7").Value = 44519
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("S7").Value = 2250
$ws.Range("D8").Value = 44505
$ws.Range("L8").Value = "Segunda"
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19500
$ws.Range("S8").Value = 2438
$ws.Range("D9").Value = 44876
$ws.Range("M9").Value = 300
$ws.Range("D10").Value = 44159
$ws.Range("N10").Value = 2000
$ws.Range("O10").Value = 2100
$ws.Range("P10").Value = 2050
$ws.Range("Q10").Value = "$/kilo (en caja de 14 kilos)"
$ws.Range("S10").Value = 2050
$ws.Range("T10").Value = 1
$ws.Range("D11").Value = 44530
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 19000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19500
$ws.Range("S11").Value = 2438
$ws.Range("D12").Value = 44530
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 16000
$ws.Range("S12").Value = 2000
$ws.Range("D13").Value = 44526
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 21000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 21000
$ws.Range("S13").Value = 2625
$ws.Range("D14").Value = 44488
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17500
$ws.Range("S14").Value = 2188
$ws.Range("D15").Value = 44880
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 22000
$ws.Range("O15").Value = 22500
$ws.Range("P15").Value = 22250
$ws.Range("S15").Value = 2781
$ws.Range("D16").Value = 44512
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 19000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 19500
$ws.Range("S16").Value = 2438
$ws.Range("D17").Value = 44509
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 19000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 19500
$ws.Range("Q17").Value = "$/bandeja 8 kilos"
$ws.Range("S17").Value = 2438
$ws.Range("T17").Value = 8
$ws.Range("D18").Value = 44894
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 22000
$ws.Range("O18").Value = 22500
$ws.Range("P18").Value = 22250
$ws.Range("S18").Value = 2781
$ws.Range("D19").Value = 44498
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 300
$ws.Range("D20").Value = 44895
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 22000
$ws.Range("O20").Value = 22500
$ws.Range("P20").Value = 22250
$ws.Range("S20").Value = 2781
$ws.Range("D21").Value = 44491
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 19000
$ws.Range("P21").Value = 18500
$ws.Range("S21").Value = 2312
$ws.Range("D22").Value = 44516
$ws.Range("M22").Value = 200
$ws.Range("O22").Value = 19000
$ws.Range("P22").Value = 18500
$ws.Range("S22").Value = 2312
$ws.Range("D23").Value = 44162
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 2000
$ws.Range("O23").Value = 2100
$ws.Range("P23").Value = 2050
$ws.Range("Q23").Value = "$/kilo (en caja de 14 kilos)"
$ws.Range("S23").Value = 2050
$ws.Range("T23").Value = 1
$ws.Range("D24").Value = 44533
$ws.Range("L24").Value = "Primera"
$ws.Range("P24").Value = 18500
$ws.Range("S24").Value = 2312
$ws.Range("D25").Value = 44533
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 16000
$ws.Range("O25").Value = 16000
$ws.Range("P25").Value = 16000
$ws.Range("S25").Value = 2000
$ws.Range("D26").Value = 44890
$ws.Range("L26").Value = "Primera"
$ws.Range("N26").Value = 22000
$ws.Range("O26").Value = 22500
$ws.Range("P26").Value = 22250
$ws.Range("S26").Value = 2781